$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the address cell from "5 rue du moulin" to "6 rue du moulin"
$ws.Range("K2").Value = "6 rue du moulin"

# Update the view: scroll so that column G becomes the left-most visible
# column (was F1) and select cell K2 (was J8)
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("K2").Select()
